$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# M4: change value from "Non" to three spaces "   "
$ws.Range("M4").Value = "   "

# J7: change date value from 43142 (2018-02-11) to 36567 (2000-02-11)
$ws.Range("J7").Value = (Get-Date -Year 2000 -Month 2 -Day 11 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# Update the active selection to J7
$ws.Range("J7").Select()
